# Timesheet updated by Kirubaharan
# Target sheet: "22-04-22" (14th worksheet / tab index 14, activeTab=13 zero-based)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(14)

# --- Row 2: add a hyperlink in a new column I, pointing to the
#     Stack Overflow question that was consulted, displaying its page title.
$ws.Rows.Item(2).RowHeight = 129.75
$i2 = $ws.Range("I2")
$i2.Value = "visual studio code - .Net SDK is not found in VSCode - Stack Overflow"
$ws.Hyperlinks.Add($i2, "https://stackoverflow.com/questions/70722884/net-sdk-is-not-found-in-vscode")

# --- Row 5: text unchanged (only shared-string index shifts because of
#     the new string above), but the row is taller in the saved layout
$ws.Rows.Item(5).RowHeight = 131.25

# --- Row 6 (Kavin): new layout / daily-update text, row grows taller
$ws.Rows.Item(6).RowHeight = 113.25
$ws.Range("B6").Value = "Layouts of TMS dashboard "
$ws.Range("C6").Value = "Layout for Login"
$ws.Range("E6").Value = "Team discussion-20 min,Meeting with Rafi-60 min,Discussion on Layout-30 min,Team discussion-30 min, Review of service-30 min"
$ws.Range("F6").Value = "Lunch and Break-90 min,Fun friday -1.5 hours, Layout Exploration -1,5 hours"

# --- Row 7 (Kirubaharan): updated daily-update text; C7 is cleared
$ws.Rows.Item(7).RowHeight = 77.25
$ws.Range("B7").Value = "Web Api exploration"
$ws.Range("C7").Value = ""
$ws.Range("E7").Value = "Team meeting - 30 mins,Meeting with Rafi - 1 hour,  Services review - 30 mins, Created Namespace ,classes and properties in VS code - 30 mins"
$ws.Range("F7").Value = "Friday activities - 1 hour 30 mins, Others (Lunch and break) - 1 hour 30 mins , Web api Tutorial - 1 hour"

# --- Row 8 (Prathima): C8 cleared; E8/F8 updated multi-line notes
$ws.Range("C8").Value = ""
$ws.Range("E8").Value = "Team Discussion(to know the progress of the team )-20mins`nReview Meeting with Rafi-1hr`nweb api (review service)-1hr`nTeam Discussion-30 mins,Created Namespace ,classes and properties in VS code - 30 mins"
$ws.Range("F8").Value = "`nLunch & Tea break - 1.5 hr`nFun session with Training team-1.5hr"

# --- Row 9 (Ragu): new registration-page layout text, row grows taller
$ws.Rows.Item(9).RowHeight = 95.25
$ws.Range("B9").Value = "registration for trainee and trainer"
$ws.Range("C9").Value = "registration for co-Ordinator"
$ws.Range("E9").Value = "team disussion  -30mim, meeting with rafi-60min,  working on layout for registration page( head)-120min, reviewed the services-35min"

# --- Row 12 (Arul): updated API/service build notes
$ws.Range("B12").Value = "Building Api for Department Service"
$ws.Range("C12").Value = "Corrected operation files for all services"
$ws.Range("E12").Value = "Team meeting -30 min , Client meeting - 60 mins ,Service operations modification - 1 hrs, Updating MOM - 15 mins"
$ws.Range("F12").Value = "lunch and others 90mins, Softskill session - 2hr"

# --- Selection moves to F7 (matching the row the author was working on)
$ws.Range("F7").Select()
